$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.763.14'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.18%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.115.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +10.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '332.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5226'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.36%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4409'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09051'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +8.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.19'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.182'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.57%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.08'
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.123.00'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +11.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.837'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.42%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.736'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.94%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '98.29'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.18%  '

$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.00%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001140'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.97%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06659'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.63%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.419'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.90%  '

$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.868.20'
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.362.17'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +11.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.261'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.22%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '163.25'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.539'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +11.53%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.94'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.09%  '

$ws.Range("E31").Value = '  +3.45%  '

$ws.Range("E32").Value = '  +2.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.250'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.08%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.537'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +28.18%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.913'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.38%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02585'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.630'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.97%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06775'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.29%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.574'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.91%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.75'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.59%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2265'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6789'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.79%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.254'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.54%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.23'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.72%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9994'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6345'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.255'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.34%  '

$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.669'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.21%  '

$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.286'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '83.37'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.45%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.59'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.33%  '
